$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 95, shifting the existing rows 95:99 down to 96:100.
$ws.Rows.Item(95).Insert()

# The data that used to be in row 95 is now in row 96; copy it into the newly
# inserted row 95 (same market/product/price data, just a new weekly date).
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(95, $c).Value = $ws.Cells.Item(96, $c).Value2
}

# Give the new row its own (one week later) date: 2021-11-09 -> 2021-11-16.
$ws.Cells.Item(95, 4).Value = 44516

# Preserve the date-number style used by column D.
$ws.Cells.Item(95, 4).NumberFormat = $ws.Cells.Item(96, 4).NumberFormat
